# "Workin on the Pie Chart"
#
# The lyric word-frequency tables on every sheet still include common
# "stop words" ("The", "And", "A") which were swamping the pie chart.
# Strip those words out of each song's word/count table. Because each
# table is already sorted descending by count, simply deleting the rows
# for the stop words leaves the remaining rows correctly ordered - no
# re-sort is needed.

$wb = $excel.ActiveWorkbook
$stopWords = @("The", "And", "A")

foreach ($ws in $wb.Worksheets) {
    $lastRow = $ws.UsedRange.Rows.Count
    for ($r = $lastRow; $r -ge 1; $r--) {
        $word = $ws.Cells.Item($r, 1).Value2
        if ($stopWords -contains $word) {
            $ws.Rows.Item($r).Delete() | Out-Null
        }
    }
}
